# Apply the latest cryptos-list snapshot (prices / 1h volume deltas) to Sheet1.
# Mirrors the GitHub Actions data refresh: most rows keep their coin/link and just
# get fresh Price (D) / Volume(1h) (E) text; rows 49-50 additionally swap which
# coin (Algorand / EnergySwap) occupies which rank.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.358.22'
$ws.Range('E2').Value = '  +2.06%  '

$ws.Range('D3').Value = '1.662.81'
$ws.Range('E3').Value = '  +1.13%  '

$ws.Range('E4').Value = '  -0.52%  '

$ws.Range('D5').Value = "'219.99"
$ws.Range('E5').Value = '  +1.30%  '

$ws.Range('D6').Value = "'0.505"
$ws.Range('E6').Value = '  +0.33%  '

$ws.Range('E7').Value = '  -0.53%  '

$ws.Range('E8').Value = '  +0.97%  '

$ws.Range('D9').Value = "'0.0627"
$ws.Range('E9').Value = '  -0.08%  '

$ws.Range('D10').Value = "'19.88"
$ws.Range('E10').Value = '  +3.67%  '

$ws.Range('E11').Value = '  +0.62%  '

$ws.Range('D12').Value = '1.897.31'
$ws.Range('E12').Value = '  +1.24%  '

$ws.Range('D13').Value = '1.663.86'
$ws.Range('E13').Value = '  +1.51%  '

$ws.Range('D14').Value = "'4.20"
$ws.Range('E14').Value = '  +0.92%  '

$ws.Range('D15').Value = "'0.534"
$ws.Range('E15').Value = '  +1.09%  '

$ws.Range('D16').Value = "'67.12"
$ws.Range('E16').Value = '  +3.91%  '

$ws.Range('D17').Value = '27.338.46'
$ws.Range('E17').Value = '  +2.04%  '

$ws.Range('D18').Value = '0.0₃0736'
$ws.Range('E18').Value = '  +0.22%  '

$ws.Range('D19').Value = "'223.81"
$ws.Range('E19').Value = '  +4.54%  '

$ws.Range('E20').Value = '  -0.59%  '

$ws.Range('D21').Value = "'6.75"
$ws.Range('E21').Value = '  +8.22%  '

$ws.Range('E22').Value = '  +1.39%  '

$ws.Range('D23').Value = "'2.52"
$ws.Range('E23').Value = '  +5.98%  '

$ws.Range('D24').Value = "'9.27"
$ws.Range('E24').Value = '  -0.39%  '

$ws.Range('D25').Value = "'147.72"
$ws.Range('E25').Value = '  +1.36%  '

$ws.Range('E26').Value = '  -0.52%  '

$ws.Range('D27').Value = "'7.41"
$ws.Range('E27').Value = '  +3.33%  '

$ws.Range('D28').Value = "'0.119"
$ws.Range('E28').Value = '  +0.92%  '

$ws.Range('D29').Value = "'16.02"
$ws.Range('E29').Value = '  +2.52%  '

$ws.Range('D30').Value = "'0.0512"
$ws.Range('E30').Value = '  +0.71%  '

$ws.Range('E32').Value = '  +1.08%  '

$ws.Range('E33').Value = '  +0.51%  '

$ws.Range('E34').Value = '  +2.27%  '

$ws.Range('D35').Value = '1.263.73'
$ws.Range('E35').Value = '  -1.68%  '

$ws.Range('E36').Value = '  +0.08%  '

$ws.Range('E37').Value = '  -0.39%  '

$ws.Range('D38').Value = "'0.537"
$ws.Range('E38').Value = '  -0.04%  '

$ws.Range('D39').Value = "'0.830"
$ws.Range('E39').Value = '  +1.51%  '

$ws.Range('E40').Value = '  -0.52%  '

$ws.Range('E41').Value = '  +0.90%  '

$ws.Range('E42').Value = '  +1.55%  '

$ws.Range('D43').Value = '1.808.75'
$ws.Range('E43').Value = '  +1.40%  '

$ws.Range('E44').Value = '  -4.53%  '

$ws.Range('E45').Value = '  +0.90%  '

$ws.Range('D46').Value = "'92.46"
$ws.Range('E46').Value = '  +0.60%  '

$ws.Range('E47').Value = '  +0.27%  '

$ws.Range('D48').Value = "'0.0517"
$ws.Range('E48').Value = '  +0.02%  '

$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = "'7.70"
$ws.Range('E49').Value = '  +0.07%  '

$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = "'0.0982"
$ws.Range('E50').Value = '  +1.49%  '

$ws.Range('D51').Value = "'0.408"
$ws.Range('E51').Value = '  +0.14%  '
